$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.887.34'
$ws.Range("E2").Value = '  +0.03%  '
$ws.Range("D3").Value = '3.437.72'
$ws.Range("E3").Value = '  -0.21%  '
$ws.Range("E4").Value = '  +0.10%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '577.94'
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '146.20'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.40%  '
$ws.Range("D7").Value = '3.437.89'
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E9").Value = '  +0.01%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '7.73'
$c.Style = "Normal"
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.123'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -1.09%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.402'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +2.16%  '
$ws.Range("D13").Value = '4.026.10'
$ws.Range("E13").Value = '  -0.26%  '
$ws.Range("E14").Value = '  +2.69%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '28.90'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -1.87%  '
$ws.Range("D16").Value = '3.438.53'
$ws.Range("E16").Value = '  -0.25%  '
$ws.Range("E17").Value = '  -1.27%  '
$ws.Range("D18").Value = '62.957.86'
$ws.Range("E18").Value = '  +0.18%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '6.34'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +1.54%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '14.33'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +0.14%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '9.17'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -1.68%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '384.84'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -2.42%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '0.558'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -0.80%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '74.34'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -1.38%  '
$ws.Range("E25").Value = '  -0.10%  '
$ws.Range("D26").Value = '3.584.29'
$ws.Range("E27").Value = '  -3.72%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '0.179'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -5.40%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '7.54'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -2.53%  '
$ws.Range("E30").Value = '  +0.23%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '8.07'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -1.60%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '2.10'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -2.14%  '
$ws.Range("E33").Value = '  -0.09%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '23.24'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -2.35%  '
$ws.Range("E35").Value = '  -9.45%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '5.27'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -1.28%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '7.05'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -0.41%  '
$ws.Range("E38").Value = '  -1.12%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '31.50'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +3.44%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '169.03'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +0.45%  '
$ws.Range("D41").Value = '3.475.65'
$ws.Range("E41").Value = '  -0.03%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.0766'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -0.35%  '
$ws.Range("E43").Value = '  -0.72%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '42.28'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -1.36%  '
$ws.Range("E45").Value = '  -1.01%  '
$ws.Range("E46").Value = '  -0.98%  '
$ws.Range("E47").Value = '  -3.23%  '
$ws.Range("D48").Value = '2.561.83'
$ws.Range("E48").Value = '  +1.78%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '2.27'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +4.02%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '6.79'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +0.73%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '22.61'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -4.42%  '
